$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''28.413.63'
$ws.Range("E2").Value = '  -0.41%  '

$ws.Range("D3").Value = '''1.823.53'
$ws.Range("E3").Value = '  -0.91%  '

$ws.Range("D4").Value = '''1.002'
$ws.Range("E4").Value = '  +0.32%  '

$ws.Range("D5").Value = '''315.04'
$ws.Range("E5").Value = '  -1.24%  '

$ws.Range("E6").Value = '  +0.30%  '

$ws.Range("D7").Value = '''0.5115'
$ws.Range("E7").Value = '  -4.16%  '

$ws.Range("D8").Value = '''0.3922'
$ws.Range("E8").Value = '  -2.29%  '

$ws.Range("D9").Value = '''0.07646'
$ws.Range("E9").Value = '  +0.64%  '

$ws.Range("D10").Value = '''41.60'
$ws.Range("E10").Value = '  -0.62%  '

$ws.Range("D11").Value = '''1.107'
$ws.Range("E11").Value = '  -0.46%  '

$ws.Range("E12").Value = '  +0.82%  '

$ws.Range("D13").Value = '''6.270'
$ws.Range("E13").Value = '  -0.89%  '

$ws.Range("D15").Value = '''7.508'
$ws.Range("E15").Value = '  -1.54%  '

$ws.Range("D16").Value = '''1.821.27'
$ws.Range("E16").Value = '  -0.46%  '

$ws.Range("D17").Value = '''93.23'
$ws.Range("E17").Value = '  +3.63%  '

$ws.Range("D18").Value = '''0.00001095'
$ws.Range("E18").Value = '  +1.97%  '

$ws.Range("D19").Value = '''0.06687'
$ws.Range("E19").Value = '  +1.32%  '

$ws.Range("D20").Value = '''17.64'
$ws.Range("E20").Value = '  -0.51%  '

$ws.Range("E21").Value = '  +0.19%  '

$ws.Range("D22").Value = '''6.143'
$ws.Range("E22").Value = '  +1.19%  '

$ws.Range("D23").Value = '''28.432.96'
$ws.Range("E23").Value = '  -0.36%  '

$ws.Range("D24").Value = '''11.18'
$ws.Range("E24").Value = '  -0.50%  '

$ws.Range("D25").Value = '''2.257'
$ws.Range("E25").Value = '  +7.31%  '

$ws.Range("D26").Value = '''20.72'
$ws.Range("E26").Value = '  +0.21%  '

$ws.Range("D27").Value = '''156.21'

$ws.Range("D28").Value = '''2.030.37'
$ws.Range("E28").Value = '  -0.75%  '

$ws.Range("E29").Value = '  -3.44%  '

$ws.Range("D30").Value = '''124.20'
$ws.Range("E30").Value = '  +0.18%  '

$ws.Range("D31").Value = '''1.108'
$ws.Range("E31").Value = '  -1.49%  '

$ws.Range("D32").Value = '''0.1088'
$ws.Range("E32").Value = '  -0.93%  '

$ws.Range("D33").Value = '''5.639'
$ws.Range("E33").Value = '  -1.28%  '

$ws.Range("D34").Value = '''3.660'
$ws.Range("E34").Value = '  -0.12%  '

$ws.Range("D35").Value = '''0.07043'
$ws.Range("E35").Value = '  -2.58%  '

$ws.Range("D36").Value = '''0.2207'
$ws.Range("E36").Value = '  -2.48%  '

$ws.Range("D37").Value = '''0.02319'
$ws.Range("E37").Value = '  -1.29%  '

$ws.Range("D38").Value = '''8.828'
$ws.Range("E38").Value = '  -0.01%  '

$ws.Range("D39").Value = '''5.156'
$ws.Range("E39").Value = '  -2.18%  '

$ws.Range("D40").Value = '''0.6242'
$ws.Range("E40").Value = '  -0.86%  '

$ws.Range("D41").Value = '''11.24'
$ws.Range("E41").Value = '  -1.20%  '

$ws.Range("D42").Value = '''1.171'
$ws.Range("E42").Value = '  -2.97%  '

$ws.Range("D43").Value = '''1.000'
$ws.Range("E43").Value = '  +0.21%  '

$ws.Range("E44").Value = '  -1.76%  '

$ws.Range("D45").Value = '''13.35'
$ws.Range("E45").Value = '  -1.49%  '

$ws.Range("D46").Value = '''0.5882'
$ws.Range("E46").Value = '  +0.47%  '

$ws.Range("D47").Value = '''3.707'
$ws.Range("E47").Value = '  -0.09%  '

$ws.Range("D48").Value = '''125.17'
$ws.Range("E48").Value = '  -0.71%  '

$ws.Range("D49").Value = '''1.976'
$ws.Range("E49").Value = '  -0.27%  '

$ws.Range("D50").Value = '''1.195'
$ws.Range("E50").Value = '  -0.11%  '

$ws.Range("D51").Value = '''0.06909'
$ws.Range("E51").Value = '  -0.33%  '
